$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.611.91'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '1.842.50'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.35'
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5272'
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3151'
$ws.Range("E8").Value = '  -3.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06807'
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.74'
$ws.Range("E10").Value = '  +0.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7811'
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '1.837.82'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007935'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '26.640.36'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").Value = '2.073.93'
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.611'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.980'
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.314'
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.29'
$ws.Range("E25").Value = '  -2.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.211'
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.677'
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.99'
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.12'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.187'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04875'
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7306'
$ws.Range("E34").Value = '  +1.81%  '
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.090'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.260'
$ws.Range("E38").Value = '  +1.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01732'
$ws.Range("E39").Value = '  -2.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4785'
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8961'
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.98'
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4161'
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.030'
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1239'
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05812'
$ws.Range("E49").Value = '  -1.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.74'
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8925'
$ws.Range("E51").Value = '  +0.73%  '
